# Auto-applied update of cryptos list values (generated from OOXML diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to remain plain text even when the value looks numeric
    # (e.g. "1.00", "253.26") so Excel does not silently convert it to a number.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '42.745.57'
Set-TextValue $ws.Range("E2") '  +3.62%  '
Set-TextValue $ws.Range("D3") '2.251.51'
Set-TextValue $ws.Range("E3") '  +3.15%  '
Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  -0.06%  '
Set-TextValue $ws.Range("D5") '253.26'
Set-TextValue $ws.Range("E5") '  -0.86%  '
Set-TextValue $ws.Range("D6") '0.636'
Set-TextValue $ws.Range("E6") '  +1.23%  '
Set-TextValue $ws.Range("D7") '70.68'
Set-TextValue $ws.Range("E7") '  +3.99%  '
Set-TextValue $ws.Range("D8") '1.00'
Set-TextValue $ws.Range("E8") '  -0.15%  '
Set-TextValue $ws.Range("D9") '0.642'
Set-TextValue $ws.Range("E9") '  +11.06%  '
Set-TextValue $ws.Range("D10") '41.30'
Set-TextValue $ws.Range("E10") '  +9.30%  '
Set-TextValue $ws.Range("D11") '59.57'
Set-TextValue $ws.Range("E11") '  +1.56%  '
Set-TextValue $ws.Range("D12") '0.0959'
Set-TextValue $ws.Range("E12") '  +2.68%  '
Set-TextValue $ws.Range("D13") '7.34'
Set-TextValue $ws.Range("E13") '  +3.01%  '
Set-TextValue $ws.Range("D14") '0.105'
Set-TextValue $ws.Range("E14") '  +1.25%  '
Set-TextValue $ws.Range("D15") '2.589.97'
Set-TextValue $ws.Range("E15") '  +3.26%  '
Set-TextValue $ws.Range("D16") '0.891'
Set-TextValue $ws.Range("E16") '  +2.32%  '
Set-TextValue $ws.Range("D17") '14.83'
Set-TextValue $ws.Range("E17") '  +2.35%  '
Set-TextValue $ws.Range("D18") '2.251.46'
Set-TextValue $ws.Range("D19") '42.739.78'
Set-TextValue $ws.Range("E19") '  +3.73%  '
Set-TextValue $ws.Range("D20") '0.0₃0978'
Set-TextValue $ws.Range("E20") '  +2.54%  '
Set-TextValue $ws.Range("D21") '6.25'
Set-TextValue $ws.Range("E21") '  +1.14%  '
Set-TextValue $ws.Range("D22") '73.04'
Set-TextValue $ws.Range("E22") '  +1.61%  '
Set-TextValue $ws.Range("D23") '235.45'
Set-TextValue $ws.Range("E23") '  +1.34%  '
Set-TextValue $ws.Range("D24") '2.12'
Set-TextValue $ws.Range("E24") '  +3.48%  '
Set-TextValue $ws.Range("D25") '3.99'
Set-TextValue $ws.Range("E25") '  +1.32%  '
Set-TextValue $ws.Range("D26") '11.63'
Set-TextValue $ws.Range("E26") '  -1.80%  '
Set-TextValue $ws.Range("E27") '  +0.16%  '
Set-TextValue $ws.Range("E28") '  -3.63%  '
Set-TextValue $ws.Range("D29") '3.68'
Set-TextValue $ws.Range("E29") '  -1.67%  '
Set-TextValue $ws.Range("E30") '  +1.58%  '
Set-TextValue $ws.Range("D31") '167.83'
Set-TextValue $ws.Range("E31") '  -0.60%  '
Set-TextValue $ws.Range("D32") '20.95'
Set-TextValue $ws.Range("E32") '  +1.56%  '
Set-TextValue $ws.Range("E33") '  +11.88%  '
Set-TextValue $ws.Range("D34") '0.124'
Set-TextValue $ws.Range("E34") '  +5.05%  '
Set-TextValue $ws.Range("E35") '  +6.33%  '
Set-TextValue $ws.Range("E36") '  +0.91%  '
Set-TextValue $ws.Range("D37") '27.95'
Set-TextValue $ws.Range("E37") '  +2.23%  '
Set-TextValue $ws.Range("D38") '4.69'
Set-TextValue $ws.Range("E38") '  +1.58%  '
Set-TextValue $ws.Range("D39") '4.19'
Set-TextValue $ws.Range("E39") '  -0.10%  '
Set-TextValue $ws.Range("D40") '0.0321'
Set-TextValue $ws.Range("E40") '  +6.81%  '
Set-TextValue $ws.Range("D41") '2.28'
Set-TextValue $ws.Range("E41") '  +3.70%  '
Set-TextValue $ws.Range("D42") '12.65'
Set-TextValue $ws.Range("E42") '  +1.13%  '
Set-TextValue $ws.Range("D43") '5.82'
Set-TextValue $ws.Range("E43") '  +2.38%  '
Set-TextValue $ws.Range("D44") '64.06'
Set-TextValue $ws.Range("E44") '  -0.69%  '
Set-TextValue $ws.Range("E45") '  -0.07%  '
Set-TextValue $ws.Range("D46") '0.203'
Set-TextValue $ws.Range("E46") '  +1.10%  '
Set-TextValue $ws.Range("D47") '8.82'
Set-TextValue $ws.Range("E47") '  +2.13%  '
Set-TextValue $ws.Range("E48") '  +1.11%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D49") '1.20'
Set-TextValue $ws.Range("E49") '  +5.20%  '
$ws.Range("B50").Value = 'BinanceUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D50") '1.00'
Set-TextValue $ws.Range("E50") '  -0.33%  '
Set-TextValue $ws.Range("D51") '4.43'
Set-TextValue $ws.Range("E51") '  +3.37%  '
